$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 5328
$ws1.Range("F12").Value = 51
$ws1.Range("F13").Value = 2219

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 5328
$ws4.Range("F15").Value = 51
$ws4.Range("F16").Value = 2219
